$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 8): Q8:T8 ---
$ws.Range("Q8").Value = "Data Type"
$ws.Range("R8").Value = "Logic"
$ws.Range("S8").Value = "Num. of Recommendations"
$ws.Range("T8").Value = "Avg. Correlation"

# --- Row 9 (Real / High) ---
$ws.Range("Q9").Value = "Real"
$ws.Range("R9").Value = "High"
$ws.Range("S9").Value = 6
$ws.Range("T9").Value = 0.79

# --- Row 10 (Random / Medium) ---
$ws.Range("Q10").Value = "Random"
$ws.Range("R10").Value = "Medium"
$ws.Range("S10").Value = 4
$ws.Range("T10").Value = 0.78

# --- Borders (thin box all around) for Q8:T10 ---
$ws.Range("Q8:T10").Borders.LineStyle = 1
$ws.Range("Q8:T10").Borders.Weight = 2

# --- Fills ---
$ws.Range("Q8:T8").Interior.Color = 8355711
$ws.Range("Q9:T9").Interior.ThemeColor = 4
$ws.Range("Q9:T9").Interior.TintAndShade = 0.59999389629810485
$ws.Range("Q10:T10").Interior.ThemeColor = 5
$ws.Range("Q10:T10").Interior.TintAndShade = 0.59999389629810485

# --- Column widths ---
$ws.Columns.Item(9).ColumnWidth = 9.28515625
$ws.Columns.Item(11).ColumnWidth = 25.7109375
$ws.Columns.Item(19).ColumnWidth = 25.7109375
$ws.Columns.Item(20).ColumnWidth = 15.42578125

# --- Sheet view ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("Q8:T10").Select()

# --- workbook path metadata ---
$wb.Path = "C:\Users\dyure\Documents\GitHub\AI_DL_Proje2\"
